$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

foreach ($ws in @($ws1, $ws4)) {
    $ws.Range("F3").Value = 1076
    $ws.Range("G4").Value = "不可售"
    $ws.Range("F6").Value = 51
    $ws.Range("F8").Value = 11164
    $ws.Range("F9").Value = 4283
    $ws.Range("F15").Value = 97
    $ws.Range("F16").Value = 14
    $ws.Range("F17").Value = 156
    $ws.Range("F18").Value = 483
    $ws.Range("F19").Value = 11222
    $ws.Range("F20").Value = 11065
}
